# Version 9.4 - Totally last fix to player ratings decoding
#
# Applies the numeric corrections to the player-ratings decoding table on
# "Sheet" (rows 4, 11, 19, 20, 25, 26, 32, 34, 39-44, 50, 52) and restores
# the workbook's on-screen selection to cell C4 (scrolled back to the
# top-left of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (column index -> new value) for every cell whose stored
# number changed in this revision.
$changes = @{
    4  = @{ 2=-560; 3=-770; 4=-700 }
    11 = @{ 2=4; 3=4; 4=4; 5=4; 6=4; 7=4; 8=4; 9=4; 10=4; 11=4; 12=4; 13=4; 14=4; 15=4; 16=4; 17=4; 18=4; 19=4; 20=4; 21=4; 22=4; 23=4; 24=4; 25=4; 26=4; 27=4; 28=4 }
    19 = @{ 2=7; 3=5; 4=5; 5=5; 6=7; 7=7; 8=7; 9=7; 10=7; 11=7; 12=7; 13=7; 14=5; 15=7; 16=7; 17=7; 18=7; 19=7; 20=5; 21=7; 22=7; 23=5; 24=7; 25=7; 26=7; 27=7; 28=5 }
    20 = @{ 2=0; 13=0 }
    25 = @{ 2=0; 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 13=0; 14=0; 15=0; 16=0; 17=0; 18=0; 19=0; 20=0; 21=0; 22=0; 23=0; 24=0; 25=0; 26=0; 27=0; 28=0 }
    26 = @{ 2=7; 3=7; 4=7; 6=7; 7=7; 13=7; 14=7; 20=7; 23=7; 26=7; 28=7 }
    32 = @{ 2=0; 13=0 }
    34 = @{ 2=0; 13=0 }
    39 = @{ 2=100; 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 14=0; 15=0; 16=0; 17=0; 18=0; 19=0; 20=0; 21=0; 22=0; 23=0; 24=0; 25=0; 26=0; 27=0; 28=0 }
    40 = @{ 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 14=0; 15=0; 16=0; 17=0; 18=0; 19=0; 20=0; 21=0; 22=0; 23=0; 24=0; 25=0; 26=0; 27=0; 28=0 }
    41 = @{ 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 14=0; 15=0; 16=0; 17=0; 18=0; 19=0; 20=0; 21=0; 22=0; 23=0; 24=0; 25=0; 26=0; 27=0; 28=0 }
    42 = @{ 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 14=0; 15=0; 16=0; 17=0; 18=0; 19=0; 20=0; 21=0; 22=0; 23=0; 24=0; 25=0; 26=0; 27=0; 28=0 }
    43 = @{ 2=0; 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 13=0; 14=0; 15=0; 16=0; 17=0; 18=0; 19=0; 20=0; 21=0; 22=0; 23=0; 24=0; 25=0; 26=0; 27=0; 28=0 }
    44 = @{ 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 14=0; 15=0; 16=0; 17=0; 18=0; 19=0; 20=0; 21=0; 22=0; 23=0; 24=0; 25=0; 26=0; 27=0; 28=0 }
    50 = @{ 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 14=0; 15=0; 16=0; 17=0; 18=0; 19=0; 20=0; 21=0; 22=0; 23=0; 24=0; 25=0; 26=0; 27=0; 28=0 }
    52 = @{ 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 14=0; 15=0; 16=0; 17=0; 18=0; 19=0; 20=0; 21=0; 22=0; 23=0; 24=0; 25=0; 26=0; 27=0; 28=0 }
}

foreach ($rowKey in $changes.Keys) {
    $colMap = $changes[$rowKey]
    foreach ($colKey in $colMap.Keys) {
        $ws.Cells.Item($rowKey, $colKey).Value = $colMap[$colKey]
    }
}

# Restore the view: select C4 (this also clears the old topLeftCell="A13"
# scroll-freeze left over from the previous save, since the new selection
# is made without any prior scrolling).
$ws.Range("C4").Select()
